$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    heading ("Play Cafelito Slot for Free - Review 2021").
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$insertionPoint = $titlePara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(2)
$newRange = $newPara.Range
$newRange.Collapse(1)  # wdCollapseStart

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Cafelito slot game and play for free. Discover the graphics, gameplay, bonuses, RTP value, and jackpots. Start playing now.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newRange.InsertXML($metaXml)

# ------------------------------------------------------------------
# 2) Remove the duplicated bold "Play Cafelito Slot for Free - Review
#    2021" paragraph that used to sit just before the meta-description
#    line near the end of the document. (Search from the end, since
#    the very first paragraph - the Heading 1 title - has the same
#    text and also renders bold via its style.)
# ------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Play Cafelito Slot for Free - Review 2021") {
        $p.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 3) Replace the old meta-description text (now italic, at the very
#    end of the document) with the new feature-image prompt text.
#    Scope the Find/Replace to the last paragraph only, since the
#    same sentence also now appears (after "Meta description: ") in
#    the paragraph inserted in step 1.
# ------------------------------------------------------------------
$oldText = "Read our review of Cafelito slot game and play for free. Discover the graphics, gameplay, bonuses, RTP value, and jackpots. Start playing now."
$newText = "Create a feature image for Cafelito featuring a happy Maya warrior with glasses in a cartoon style. The background should have a coffee shop theme with images of coffee beans, cups, and machines. The Maya warrior should be holding a cup of coffee with a smile on their face. They should be wearing a colorful outfit with traditional Maya patterns, and their hair should be decorated with coffee beans and flowers. The image should be bright and vibrant, capturing the fun and excitement of playing the Cafelito slot game."

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
